# Adds the Arabic pangram (transliteration + Arabic script) as two new
# lines below the existing translation-test table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A18").Value = "dahab al-thaalib al-bunni sarie abr al-ghaba lajalab shaheb minn al-maa "
$ws.Range("A19").Value = "ذهب الثعلب البني السريع عبر الغابة لجلب شاحب من الماء"

$ws.Range("A18").WrapText = $true
$ws.Rows.Item(18).RowHeight = 24.05

[void]$ws.Range("C18").Select()
